# Add a new dictionary entry (key / de / en) for YOU_FINISHED, to enforce
# clearing of the timeout for SRS and ART ("clear_page").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = "YOU_FINISHED"
$ws.Range("B36").Value = "Sie haben den {{test_name}} beendet."
$ws.Range("C36").Value = "You finished the {{test_name}}."

# Scroll the view down so row 16 is at the top, and leave the active
# selection on B36, matching the author's final view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("B36").Select()
